$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (before) values for rows 2, 3, 4 for the columns that change.
$cols = @("A","B","E","F","G","H","Q","R","AI","AO")

$rowData = @{}
foreach ($r in 2..4) {
    $data = @{}
    foreach ($col in $cols) {
        $data[$col] = $ws.Range("$col$r").Value2
    }
    $rowData[$r] = $data
}

# New row 2 gets old row 3's data; new row 3 gets old row 4's data; new row 4 gets old row 2's data.
$mapping = @{ 2 = 3; 3 = 4; 4 = 2 }

foreach ($destRow in 2..4) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowData[$srcRow][$col]
    }
}
